# SnakeChaosHouse Members.xlsx - add "Password" column (R) for every member row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell R1: "Password" with bold font + thin box border + centered/top aligned ---
$r1 = $ws.Range("R1")
$r1.Value = "Password"
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108   # xlCenter
$r1.VerticalAlignment = -4160     # xlTop
$r1.Borders.LineStyle = 1         # xlContinuous
$r1.Borders.Weight = 2            # xlThin

# --- Data rows: one password per existing member row ---
# Rows that already carried an (empty) styled placeholder cell in column R need
# their style reset back to Normal/default before the value is written, to match
# the plain (unstyled) password cells produced by the edit.
$ws.Range("R2").Style = "Normal"
$ws.Range("R2").Value = "D85mk4Kh"
$ws.Range("R3").Style = "Normal"
$ws.Range("R3").Value = "CdmYUljS"
$ws.Range("R4").Style = "Normal"
$ws.Range("R4").Value = "FtXawwUW"
$ws.Range("R5").Style = "Normal"
$ws.Range("R5").Value = "AkiI7rOW"
$ws.Range("R6").Style = "Normal"
$ws.Range("R6").Value = "mVbMC7n2"
$ws.Range("R7").Style = "Normal"
$ws.Range("R7").Value = "tPypM1Ne"
$ws.Range("R8").Style = "Normal"
$ws.Range("R8").Value = "BSJVo8jd"
$ws.Range("R9").Style = "Normal"
$ws.Range("R9").Value = "DEIoeIfu"
$ws.Range("R10").Style = "Normal"
$ws.Range("R10").Value = "qI2pgPyk"
$ws.Range("R11").Style = "Normal"
$ws.Range("R11").Value = "DUMfHEKG"
$ws.Range("R12").Style = "Normal"
$ws.Range("R12").Value = "zoqJPTKD"
$ws.Range("R13").Style = "Normal"
$ws.Range("R13").Value = "5EYMDOFM"
$ws.Range("R14").Style = "Normal"
$ws.Range("R14").Value = "HLqiCHBU"
$ws.Range("R18").Value = "sqm5xGG8"
$ws.Range("R19").Value = "YHNvgkDg"
$ws.Range("R20").Value = "UFNfsLp4"
$ws.Range("R21").Value = "O8v66rwM"
$ws.Range("R22").Value = "yrSfSCq7"
$ws.Range("R23").Value = "LRRjkObr"
$ws.Range("R27").Value = "Y7BMTHUN"
$ws.Range("R28").Value = "aiKJySsa"
$ws.Range("R29").Value = "aj4LoLfL"
$ws.Range("R30").Value = "h85MBOW5"
$ws.Range("R31").Value = "mTdNscvN"
$ws.Range("R32").Value = "j3ya7T5o"
$ws.Range("R33").Value = "YqGrdV0T"
$ws.Range("R34").Value = "DoNKpVL5"
$ws.Range("R35").Value = "NfLz5umL"
$ws.Range("R36").Value = "Zrv5gpgj"
$ws.Range("R37").Value = "sTWiKbV7"
$ws.Range("R38").Value = "s1d1uTLw"
$ws.Range("R39").Value = "xH30UiQ5"
$ws.Range("R40").Value = "124k6x3A"
$ws.Range("R41").Value = "pBs9E69D"
$ws.Range("R42").Value = "7G4PLLOt"
$ws.Range("R43").Value = "AKddEp87"
$ws.Range("R44").Style = "Normal"
$ws.Range("R44").Value = "m4gKPwGG"
$ws.Range("R45").Value = "Zze6WFWe"
$ws.Range("R49").Value = "QxKLma6O"
$ws.Range("R50").Value = "J4LH3d7N"
$ws.Range("R51").Value = "hHjirkz8"
$ws.Range("R52").Value = "Jz1XKzPO"
$ws.Range("R53").Value = "YxCI3niz"
$ws.Range("R54").Value = "CO9dP7Ts"
$ws.Range("R55").Value = "JwswhoBU"
$ws.Range("R56").Value = "vGidvj1X"
$ws.Range("R57").Value = "XIU89l07"
$ws.Range("R58").Value = "akKd9Wt6"
$ws.Range("R59").Value = "jvKwjwrF"
$ws.Range("R60").Value = "SdPYZlg7"
$ws.Range("R61").Value = "x9b6YYG1"
$ws.Range("R62").Value = "BDiwelI7"
$ws.Range("R63").Value = "pf5IhhH4"
$ws.Range("R64").Value = "7c9ZJ5CZ"
$ws.Range("R65").Value = "FmvCMhWX"
$ws.Range("R66").Value = "wxgkdl8v"

# --- Column R width ---
$ws.Columns.Item(18).ColumnWidth = 23 + 1/6

# --- Selection / view state ---
$ws.Range("G33").Select()
